# "moved project direct sort" - insert a new "Odd/Even sort" results column
# before the existing "Direct Sort" timing column (which shifts from G to H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at G; this pushes the old column G (Direct Sort
# timing header, bold) to column H, including its style/value and shifting
# the dimension/used range accordingly.
$null = $ws.Columns("G:G").Insert()

# Fill in the header for the newly inserted column.
$ws.Range("G1").Value = "Odd/Even sort"

# Give the new column the same width as column F (both end up 28.33203125
# wide in the original file).
$ws.Columns("G:G").ColumnWidth = $ws.Columns("F:F").ColumnWidth()

# Move the active selection to the newly added cell.
$null = $ws.Range("G1").Select()
